$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells AN1:AX1 - copy the bold/bordered header style from an
#    existing header cell (A1) first, then set the values.
# ---------------------------------------------------------------------------
$headers = @(
    "Model Type",
    "Timeout (s)",
    "Temperature",
    "Max Tokens",
    "Top P",
    "Top K",
    "Frequency Penalty",
    "Presence Penalty",
    "N",
    "Num Ctx",
    "Format"
)

$headerStartCol = 40  # column AN
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $headerStartCol + $i
    $cell = $ws.Cells.Item(1, $col)
    $ws.Range("A1").Copy($cell)
    $cell.Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2. New data row 27 - fully populated, including the new AN:AX columns.
# ---------------------------------------------------------------------------
$row27 = @(
    "mistral:7b-instruct-v0.3-q5_K_M",
    "llama3:70b",
    10,
    200,
    50.4,
    4.41,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_match.txt",
    4.41,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_correct.txt",
    4.41,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_executable.txt",
    0,
    0,
    6.23,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_match.txt",
    6.23,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_correct.txt",
    6.23,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_executable.txt",
    0,
    31.03,
    8.72,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_match.txt",
    8.72,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_correct.txt",
    8.72,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_executable.txt",
    0,
    2,
    2,
    "text",
    120,
    0.2,
    300,
    1,
    10,
    1,
    1.2,
    1,
    1024,
    "json"
)

# ---------------------------------------------------------------------------
# 3. New data row 28 - populated through AM only (AN:AX left blank, as in
#    the source diff).
# ---------------------------------------------------------------------------
$row28 = @(
    "mistral:7b-instruct-v0.3-q5_K_M",
    "llama3:70b",
    10,
    200,
    51.26,
    5.25,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_match_2.txt",
    5.25,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_correct_2.txt",
    5.25,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_executable_2.txt",
    0,
    0,
    5.13,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_match_2.txt",
    5.13,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_correct_2.txt",
    5.13,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_executable_2.txt",
    0,
    32.64,
    8.22,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_match_2.txt",
    8.22,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_correct_2.txt",
    8.22,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_executable_2.txt",
    0,
    2,
    2
)

# ---------------------------------------------------------------------------
# 4. New data row 29 - populated through AM only (AN:AX left blank, as in
#    the source diff).
# ---------------------------------------------------------------------------
$row29 = @(
    "mistral:7b-instruct-v0.3-q5_K_M",
    "llama3:70b",
    10,
    200,
    333.52,
    14.25,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_match_3.txt",
    14.25,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_correct_3.txt",
    14.25,
    0.5,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_executable_3.txt",
    0.8571428571428571,
    0,
    14.7,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_match_3.txt",
    14.7,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_correct_3.txt",
    14.7,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_fewshot_executable_3.txt",
    1,
    285.87,
    18.68,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_match_3.txt",
    18.68,
    1,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_correct_3.txt",
    18.68,
    0,
    "logs\mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_10_200_test_bootstrap_executable_3.txt",
    0.8571428571428571,
    2,
    2
)

function Write-RowValues($ws, $rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

Write-RowValues $ws 27 $row27
Write-RowValues $ws 28 $row28
Write-RowValues $ws 29 $row29

Write-Host "Applied log_evaluations edits: added AN:AX header columns and rows 27-29"
